# Update cryptos list: refresh Price (D) and Volume(1h) (E) values;
# rows 35/36 also swap ARBITRUM/WEMIXToken (Coin name + Link).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.383.26"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "'2.612.88"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'308.24"
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("D6").Value = "'100.42"
$ws.Range("E6").Value = "  +4.00%  "
$ws.Range("D7").Value = "'0.603"
$ws.Range("E7").Value = "  +3.26%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.580"
$ws.Range("E9").Value = "  +7.26%  "
$ws.Range("D10").Value = "'39.63"
$ws.Range("E10").Value = "  +8.31%  "
$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = "  +6.17%  "
$ws.Range("D12").Value = "'54.20"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "'8.18"
$ws.Range("E13").Value = "  +8.35%  "
$ws.Range("D14").Value = "'3.016.09"
$ws.Range("E14").Value = "  +3.79%  "
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "'2.619.52"
$ws.Range("E16").Value = "  +4.13%  "
$ws.Range("D17").Value = "'0.922"
$ws.Range("E17").Value = "  +4.84%  "
$ws.Range("D18").Value = "'14.98"
$ws.Range("E18").Value = "  +3.27%  "
$ws.Range("D19").Value = "'46.549.62"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("E20").Value = "  +5.33%  "
$ws.Range("D21").Value = "'12.93"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("E22").Value = "  +4.08%  "
$ws.Range("D23").Value = "'71.60"
$ws.Range("E23").Value = "  +4.58%  "
$ws.Range("D24").Value = "'274.67"
$ws.Range("E24").Value = "  +10.10%  "
$ws.Range("E25").Value = "  +7.00%  "
$ws.Range("E26").Value = "  +6.90%  "
$ws.Range("D27").Value = "'29.03"
$ws.Range("E27").Value = "  +28.53%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'4.02"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").Value = "'10.63"
$ws.Range("E30").Value = "  +5.85%  "
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("D32").Value = "'39.11"
$ws.Range("E32").Value = "  -3.53%  "
$ws.Range("D33").Value = "'6.38"
$ws.Range("E33").Value = "  +11.27%  "
$ws.Range("D34").Value = "'3.67"
$ws.Range("E34").Value = "  -5.11%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'2.25"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.86"
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("E37").Value = "  +6.00%  "
$ws.Range("D38").Value = "'151.07"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("E39").Value = "  +5.46%  "
$ws.Range("D40").Value = "'0.124"
$ws.Range("E40").Value = "  +5.31%  "
$ws.Range("D41").Value = "'23.48"
$ws.Range("E41").Value = "  +40.82%  "
$ws.Range("D42").Value = "'15.89"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").Value = "'3.67"
$ws.Range("E43").Value = "  +9.43%  "
$ws.Range("D44").Value = "'0.0333"
$ws.Range("E44").Value = "  +8.06%  "
$ws.Range("D45").Value = "'4.08"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "'2.127.22"
$ws.Range("E46").Value = "  +5.22%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "'93.61"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").Value = "'9.53"
$ws.Range("E49").Value = "  +8.57%  "
$ws.Range("D50").Value = "'1.78"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "'109.31"
$ws.Range("E51").Value = "  +3.05%  "
